$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "лошок"
$ws.Range("C4").Value = 64
